$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pin Map")

# --- Move the "not connected" annotations from rows 51/52/55/56 up to rows 27/28/33/34 ---

# Row 27 <- data that was on row 51 (_IRQ / In/Out)
$ws.Cells.Item(27, 10).Value = "_IRQ"
$ws.Cells.Item(27, 11).Value = "In/Out"
$ws.Cells.Item(27, 12).Value = "Not Connected/unused"
$ws.Cells.Item(27, 12).Interior.Color = 65535

# Row 28 <- data that was on row 52 (_NMI / Out)
$ws.Cells.Item(28, 10).Value = "_NMI"
$ws.Cells.Item(28, 11).Value = "Out"
$ws.Cells.Item(28, 12).Value = "Not Connected/unused"
$ws.Cells.Item(28, 12).Interior.Color = 65535

# Row 33 <- data that was on row 55 (_DMA / Out)
$ws.Cells.Item(33, 10).Value = "_DMA"
$ws.Cells.Item(33, 11).Value = "Out"
$ws.Cells.Item(33, 12).Value = "Not Connected/unused"
$ws.Cells.Item(33, 12).Interior.Color = 65535

# Row 34 <- data that was on row 56 (DOT Clk / In)
# J34 did not exist as a cell before, so give it the plain bordered style too.
$ws.Cells.Item(34, 10).Value = "DOT Clk"
$ws.Cells.Item(34, 13).Copy()
$ws.Cells.Item(34, 10).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(34, 11).Value = "In"
$ws.Cells.Item(34, 12).Value = "Not Connected/unused"
$ws.Cells.Item(34, 12).Interior.Color = 65535

# Clear out the old entries on rows 51/52/55/56 (value + format reset to plain bordered cell)
foreach ($r in 51, 52, 55, 56) {
    $ws.Cells.Item($r, 10).Value = $null
    $ws.Cells.Item($r, 11).Value = $null
    $ws.Cells.Item($r, 12).Value = $null
    $ws.Cells.Item(53, 13).Copy()
    $ws.Cells.Item($r, 12).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Mark the SD-card pins (now routed through flash/SRAM instead) ---
foreach ($r in 45..50) {
    $ws.Cells.Item($r, 8).Interior.Color = 49407
    $ws.Cells.Item($r, 14).Value = "SD Card"
}

# Rows 51-57: tag column N with "Flash/SRAM locations"
foreach ($r in 51..57) {
    $ws.Cells.Item($r, 14).Value = "Flash/SRAM locations"
}

# --- Selection bookkeeping ---
$ws.Range("I9").Select()
